# Apply the diff: for rows 35-68 in worksheet "Exp_data", move the
# "K_LAENDER" / "A_LAENDER_xx" pair that currently lives in columns E/F
# into columns C/D (overwriting the "K_SEX" / "A_SEX_x" values that were
# there), and blank out columns E/F. For the two rows of each block that
# had no K_LAENDER/A_LAENDER pair (the first row of each K_SEX group,
# i.e. rows 35 and 52), simply blank out C/D as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Exp_data")

for ($row = 35; $row -le 68; $row++) {
    $cC = $ws.Cells.Item($row, 3)   # column C
    $cD = $ws.Cells.Item($row, 4)   # column D
    $cE = $ws.Cells.Item($row, 5)   # column E
    $cF = $ws.Cells.Item($row, 6)   # column F

    $eVal = $cE.Value2
    $fVal = $cF.Value2

    if ($eVal -eq $null) { $eVal = "" }
    if ($fVal -eq $null) { $fVal = "" }

    # Move E/F contents into C/D, then clear E/F.
    $cC.Value = $eVal
    $cD.Value = $fVal
    $cE.Value = ""
    $cF.Value = ""
}
